$wb = $excel.ActiveWorkbook

# Delete Sheet2 and Sheet3, leaving only Sheet1
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Sheet2").Delete()
[void]$wb.Worksheets.Item("Sheet3").Delete()
$excel.DisplayAlerts = $true

# Rename Sheet1 to "Sales Split %"
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Sales Split %"
